$wb = $excel.ActiveWorkbook

# The "Metadata" sheet holds the Property/Value table
$ws = $wb.Worksheets.Item("Metadata")

# Row 5 is "Title" - it was empty, set it to the same value as "Name" (LangueParlee)
$ws.Range("B5").Value = "LangueParlee"

# Row 8 is "Date" - bump the generation timestamp
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
